$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Radiology", "Unlikely", "Doctors are patronising and made me feel bad"),
    @("Theatre Treatment Suite Implants", "Extremely Unlikely", "doctors dont seem to care about me, felt ignored"),
    @("Labour and Delivery Suite", "Extremely Unlikely", "Clenliness isn't the best but otherwise okay"),
    @("Labour and Delivery Suite", "Extremely Unlikely", "doctors dont seem to care about me, felt ignored"),
    @("Gynaecology", "Extremely Unlikely", "Service recieved was adaquete but staff seemed like they dont care at all"),
    @("Theatre Treatment Suite Implants", "Unlikely", "Waited over 5 hours"),
    @("Radiology", "Unlikely", "doctors dont seem to care about me, felt ignored"),
    @("Theatre Treatment Suite Implants", "Extremely Unlikely", "Waited too long to find a parking spot"),
    @("Special Care Baby Unit", "Extremely Unlikely", "Food was terrible"),
    @("Theatre Treatment Suite Implants", "Unlikely", "I waited for a long time to be turned away"),
    @("Special Care Baby Unit", "Extremely Unlikely", "Waited for long time for poor service"),
    @("A&E", "Unlikely", "Felt as if i was not a priority"),
    @("Day Surgery", "Extremely Unlikely", "Service recieved was adaquete but staff seemed like they dont care at all"),
    @("Rehab Services", "Unlikely", "Waited over 5 hours"),
    @("Special Care Baby Unit", "Extremely Unlikely", "Service recieved was adaquete but staff seemed like they dont care at all"),
    @("Sitwell", "Unlikely", "doctors dont seem to care about me, felt ignored"),
    @("Radiology", "Extremely Unlikely", "Long wait times")
)

$row = 2
foreach ($r in $data) {
    $ws.Cells.Item($row, 1).Value = $r[0]
    $ws.Cells.Item($row, 2).Value = $r[1]
    $ws.Cells.Item($row, 3).Value = $r[2]
    $row = $row + 1
}

$ws.Range("A19:XFD101").Select()
